# Update the "Metadata" sheet (StructureDefinition metadata for
# ror-territorial-division) to the new release values.
$wb = $excel.ActiveWorkbook

$meta = $wb.Worksheets.Item("Metadata")
$meta.Range("B3").Value  = "0.4.0-snapshot-1"                   # Version
$meta.Range("B6").Value  = "draft"                               # Status
$meta.Range("B8").Value  = "2024-05-23T12:16:26+00:00"           # Date
$meta.Range("B10").Value = "ANS (https://esante.gouv.fr)"        # Contact

# On the "Elements" sheet, the two "Mapping" columns (AK = column 37,
# AL = column 38) are swapped: the business-mapping column moves in
# front of the RIM-mapping column. Swap the header, every data row,
# and the column widths together so the whole columns trade places.
$elements = $wb.Worksheets.Item("Elements")

$lastRow = 16
For ($r = 1; $r -le $lastRow; $r++) {
    $akCell = $elements.Cells.Item($r, 37)
    $alCell = $elements.Cells.Item($r, 38)
    $akVal = $akCell.Value2
    $alVal = $alCell.Value2
    $akCell.Value = $alVal
    $alCell.Value = $akVal
}

# Column widths also swap along with the data (AK becomes the wide
# business-mapping column, AL becomes the narrower RIM-mapping column).
$elements.Range("AK1").ColumnWidth = 73.05
$elements.Range("AL1").ColumnWidth = 24.15
